# Add "Mismatch Tolerance (Levenshtein Distance)" option columns to the config sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: new header label in J1
$ws.Range("J1").Value = "Mismatch Tolerance (Levenshtein Distance)"

# Row 2: sub-headers for the new columns
$ws.Range("J2").Value = "Leader Proximal"
$ws.Range("K2").Value = "Leader Distal"

# Row 3: shift the old example text from J3 to L3, and put numeric tolerance
# values (1) in the two new columns J3 and K3.
$ws.Range("L3").Value = $ws.Range("J3").Value2
$ws.Range("L3").Font.Italic = $ws.Range("J3").Font.Italic

$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("J3").Font.Italic = $true
$ws.Range("K3").Font.Italic = $true

# Update the selection to match the saved view state.
$ws.Range("G11").Select()
